$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add a new row of piping input data (row 4), duplicating row 3 but with
# an incremented "North" value (Line no. stays "100-test").
$ws.Range("A4").Value = "100-test"
$ws.Range("B4").Value = 101
$ws.Range("C4").Value = 100
$ws.Range("D4").Value = 10
$ws.Range("E4").Value = 10
$ws.Range("F4").Value = 28
$ws.Range("G4").Value = 0
$ws.Range("H4").Value = 0.3

# Move the active selection to C4, matching the author's final cursor position.
$ws.Range("C4").Select()
